# "have counties working now"
# Appends two new days (rows 19 & 20) of COVID case data to Sheet1, extending
# the running tables/shared formulas down through row 20, and updates the
# window selection to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows right after the current last row (18) so that the
# row-18 number formats/styles (date format in col A, percent format in
# C/F/Y/Z) carry down automatically, exactly like Excel does when a user
# types a new row directly under an existing table.
$ws.Rows("19:20").Insert()

# --- Row 19 (2020-03-22, serial 43922) ---
$ws.Range("A19").Value = 43922
$ws.Range("B19").Value = 7738
$ws.Range("C19").Formula = "=(B19-B18)/B18"
$ws.Range("D19").Value = 122
$ws.Range("E19").Value = 682
$ws.Range("F19").Formula = "=E19/B19"
$ws.Range("U19").Value = 51738
$ws.Range("V19").Formula = "=B19"
$ws.Range("W19").Formula = "=U19-U18"
$ws.Range("X19").Formula = "=V19-V18"
$ws.Range("Y19").Formula = "=X19/W19"
$ws.Range("Z19").Formula = "=V19/U19"

# --- Row 20 (2020-03-23, serial 43923) ---
$ws.Range("A20").Formula = "=A19+1"
$ws.Range("B20").Value = 8966
$ws.Range("C20").Formula = "=(B20-B19)/B19"
$ws.Range("D20").Value = 154
$ws.Range("E20").Value = 813
$ws.Range("F20").Formula = "=E20/B20"
$ws.Range("U20").Value = 56608
$ws.Range("V20").Formula = "=B20"
$ws.Range("W20").Formula = "=U20-U19"
$ws.Range("X20").Formula = "=V20-V19"
$ws.Range("Y20").Formula = "=X20/W20"
$ws.Range("Z20").Formula = "=V20/U20"

# Recalculate everything so cached <v> results are fresh.
$excel.Calculate()

# Restore the cursor/selection to where the author left it.
$ws.Range("W26").Select()
